# The row-1 header of this import-template sheet lists one shared string
# per generated column: protocol, lbl, is_locked_lbl, is_enabled_lbl,
# order_by, rem (columns A-F). This change drops the is_locked_lbl and
# is_enabled_lbl columns (C and D) entirely, so order_by/rem shift left
# into C1/D1 and the sheet ends up with only 4 used columns (A-D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftToLeft = -4159: remove C1:D1 and slide the remaining cells of the
# row (order_by, rem) left to fill the gap, instead of just blanking them.
$ws.Range("C1:D1").Delete(-4159)
